# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to match the data refresh captured in the commit.

$wb = $excel.ActiveWorkbook

# Changes common to both "展览" and "全部类型" sheets (row -> new F value)
$commonChanges = @{
    3  = 3368
    4  = 248
    6  = 216
    7  = 1732
    8  = 1656
    9  = 476
    14 = 38
    16 = 3
    19 = 235
    23 = 62
    24 = 40
    25 = 24
    26 = 401
    27 = 261
    31 = 30
    32 = 428
    33 = 2295
    36 = 482
    37 = 562
    38 = 566
    39 = 432
    40 = 235
    42 = 417
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonChanges.Keys) {
        $ws.Range("F$row").Value = $commonChanges[$row]
    }

    # F17 differs between the two sheets before the edit, but both become 29
    $ws.Range("F17").Value = 29
}

$wb.Save()
